$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 42-49
$ws.Range("A42:A49").Value = "S6"
$ws.Range("B42:B49").Value = "Yelena"
$ws.Range("C42:C49").Value = "18/7/2013"
$ws.Range("D42:D49").Value = "2013-07-18-yelena"
$ws.Range("E42").Value = "2013-07-18-15-11-12"
$ws.Range("E43").Value = "2013-07-18-15-18-53"
$ws.Range("E44").Value = "2013-07-18-15-26-03"
$ws.Range("E45").Value = "2013-07-18-15-32-16"
$ws.Range("E46").Value = "2013-07-18-15-43-31"
$ws.Range("E47").Value = "2013-07-18-15-51-54"
$ws.Range("E48").Value = "2013-07-18-15-58-42"
$ws.Range("E49").Value = "2013-07-18-16-05-13"
$ws.Range("F42").Value = "ssvep-15Hz"
$ws.Range("G42").Value = 15
$ws.Range("H42").Value = 0
$ws.Range("F43").Value = "hybrid-8-57Hz"
$ws.Range("G43").Value = 8.57
$ws.Range("H43").Value = 1
$ws.Range("F44").Value = "hybrid-10Hz"
$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 1
$ws.Range("F45").Value = "ssvep-10Hz"
$ws.Range("G45").Value = 10
$ws.Range("H45").Value = 0
$ws.Range("F46").Value = "ssvep-12Hz"
$ws.Range("G46").Value = 12
$ws.Range("H46").Value = 0
$ws.Range("F47").Value = "hybrid-15Hz"
$ws.Range("G47").Value = 15
$ws.Range("H47").Value = 1
$ws.Range("F48").Value = "ssvep-8-57Hz"
$ws.Range("G48").Value = 8.57
$ws.Range("H48").Value = 0
$ws.Range("F49").Value = "hybrid-12Hz"
$ws.Range("G49").Value = 12
$ws.Range("H49").Value = 1

# Rows 50-57
$ws.Range("A50:A57").Value = "S7"
$ws.Range("B50:B57").Value = "Robert"
$ws.Range("C50:C57").Value = "30/7/2013"
$ws.Range("D50:D57").Value = "2013-07-30-robert"
$ws.Range("E50").Value = "2013-07-30-14-33-23"
$ws.Range("E51").Value = "2013-07-30-14-41-12"
$ws.Range("E52").Value = "2013-07-30-15-01-18"
$ws.Range("E53").Value = "2013-07-30-15-08-27"
$ws.Range("E54").Value = "2013-07-30-15-20-51"
$ws.Range("E55").Value = "2013-07-30-15-28-11"
$ws.Range("E56").Value = "2013-07-30-15-40-15"
$ws.Range("E57").Value = "2013-07-30-15-47-39"
$ws.Range("F50").Value = "hybrid-15Hz"
$ws.Range("G50").Value = 15
$ws.Range("H50").Value = 1
$ws.Range("F51").Value = "ssvep-15Hz"
$ws.Range("G51").Value = 15
$ws.Range("H51").Value = 0
$ws.Range("F52").Value = "hybrid-8-57Hz"
$ws.Range("G52").Value = 8.57
$ws.Range("H52").Value = 1
$ws.Range("F53").Value = "hybrid-12Hz"
$ws.Range("G53").Value = 12
$ws.Range("H53").Value = 1
$ws.Range("F54").Value = "ssvep-8-57Hz"
$ws.Range("G54").Value = 8.57
$ws.Range("H54").Value = 0
$ws.Range("F55").Value = "ssvep-10Hz"
$ws.Range("G55").Value = 10
$ws.Range("H55").Value = 0
$ws.Range("F56").Value = "hybrid-10Hz"
$ws.Range("G56").Value = 10
$ws.Range("H56").Value = 1
$ws.Range("F57").Value = "ssvep-12Hz"
$ws.Range("G57").Value = 12
$ws.Range("H57").Value = 0

# Rows 58-65
$ws.Range("A58:A65").Value = "S8"
$ws.Range("B58:B65").Value = "Alejandro"
$ws.Range("E58").Value = "2013-08-07-10-48-11"
$ws.Range("E59").Value = "2013-08-07-10-58-05"
$ws.Range("E60").Value = "2013-08-07-11-04-37"
$ws.Range("E61").Value = "2013-08-07-11-11-34"
$ws.Range("E62").Value = "2013-08-07-11-26-57"
$ws.Range("E63").Value = "2013-08-07-11-33-46"
$ws.Range("E64").Value = "2013-08-07-11-40-19"
$ws.Range("E65").Value = "2013-08-07-11-47-18"
$ws.Range("D58:D65").Value = "2013-08-07-alejandro"
$ws.Range("F58").Value = "hybrid-8-57Hz"
$ws.Range("G58").Value = 8.57
$ws.Range("H58").Value = 1
$ws.Range("F59").Value = "ssvep-15Hz"
$ws.Range("G59").Value = 15
$ws.Range("H59").Value = 0
$ws.Range("F60").Value = "ssvep-12Hz"
$ws.Range("G60").Value = 12
$ws.Range("H60").Value = 0
$ws.Range("F61").Value = "hybrid-15Hz"
$ws.Range("G61").Value = 15
$ws.Range("H61").Value = 1
$ws.Range("F62").Value = "hybrid-12Hz"
$ws.Range("G62").Value = 12
$ws.Range("H62").Value = 1
$ws.Range("F63").Value = "ssvep-8-57Hz"
$ws.Range("G63").Value = 8.57
$ws.Range("H63").Value = 0
$ws.Range("F64").Value = "hybrid-10Hz"
$ws.Range("G64").Value = 10
$ws.Range("H64").Value = 1
$ws.Range("F65").Value = "ssvep-10Hz"
$ws.Range("G65").Value = 10
$ws.Range("H65").Value = 0

# Fix C58:C65 to be numeric dates with the same style as C34 (m/d/yyyy)
$ws.Range("C58:C65").Value = 41463
$ws.Range("C34").Copy()
$ws.Range("C58:C65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match final state
$ws.Range("A59:D65").Select()
$excel.ActiveWindow.ScrollRow = 40
